$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 51 (brand new row): copy formatting for the whole row from row 49,
#     which is the most recent row sharing the same "6.4.2020" date block ---
$ws.Range("A49:G49").Copy($ws.Range("A51")) | Out-Null

# --- Row 52 already existed (it only held a lone, empty D52 cell styled
#     s="5"). Copy formatting for A:C and E:G from row 49, but deliberately
#     leave the D column alone so D52 keeps its original s="5" styling
#     instead of picking up row 49's s="3". ---
$ws.Range("A49:C49").Copy($ws.Range("A52")) | Out-Null
$ws.Range("E49:G49").Copy($ws.Range("E52")) | Out-Null

# Values for row 51 (A51 keeps the "6.4.2020" shared string inherited from
# the row-49 copy, so it is left untouched here)
$ws.Range("B51").Value2 = 0.51388888888888895
$ws.Range("C51").Value2 = 0.52777777777777779
$ws.Range("D51").Formula = "=C51-B51"
$ws.Range("E51").Value2 = "Pattern Gen 1"
$ws.Range("F51").Value2 = "Improve Code"
$ws.Range("G51").Value2 = "Add packages"

# Values for row 52 (A52 keeps the "6.4.2020" shared string inherited from
# the row-49 copy, so it is left untouched here)
$ws.Range("B52").Value2 = 0.55555555555555558
$ws.Range("C52").Value2 = 0.5625
$ws.Range("D52").Formula = "=C52-B52"
$ws.Range("E52").Value2 = "Pattern Gen 2"
$ws.Range("F52").Value2 = "Imrpove Code"
$ws.Range("G52").Value2 = "Add packages"

# Move / record the active selection as it ends up after the edit.
$ws.Range("G53").Select() | Out-Null

Write-Host "edit applied"
